# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Repollo" (Feria Lagunitas de Puerto
# Montt) right before the current row 355, pushing the existing rows
# 355-375 down to 357-377.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 355 (each Insert() call shifts
# everything at/after row 355 down by one row, inheriting the formatting
# of the row above - same behaviour as native Excel "Insert Sheet Rows").
$ws.Rows(355).Insert()
$ws.Rows(355).Insert()

# --- New row 355 ---------------------------------------------------------
$ws.Range("A355").Value = 4
$ws.Range("B355").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C355").Value = "Los Lagos"
$ws.Range("D355").Value = 44610
$ws.Range("E355").Value = 10
$ws.Range("F355").Value = 100112006
$ws.Range("G355").Value = "Repollo"
$ws.Range("H355").Value = "Copenhague"
$ws.Range("I355").Value = "Primera"
$ws.Range("J355").Value = 400
$ws.Range("K355").Value = 2000
$ws.Range("L355").Value = 2000
$ws.Range("M355").Value = 2000
$ws.Range("N355").Value = "$/unidad"
$ws.Range("O355").Value = "Región Metropolitana"
$ws.Range("P355").Value = 2000
$ws.Range("Q355").Value = 1
$ws.Range("R355").Value = "Hortaliza"

# --- New row 356 ---------------------------------------------------------
$ws.Range("A356").Value = 4
$ws.Range("B356").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C356").Value = "Los Lagos"
$ws.Range("D356").Value = 44610
$ws.Range("E356").Value = 10
$ws.Range("F356").Value = 100112006
$ws.Range("G356").Value = "Repollo"
$ws.Range("H356").Value = "Crespo record"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 600
$ws.Range("K356").Value = 1700
$ws.Range("L356").Value = 1700
$ws.Range("M356").Value = 1700
$ws.Range("N356").Value = "$/unidad"
$ws.Range("O356").Value = "Región Metropolitana"
$ws.Range("P356").Value = 1700
$ws.Range("Q356").Value = 1
$ws.Range("R356").Value = "Hortaliza"
